$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header row (copy format from H1 so the same style index is reused)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data for columns I and J, rows 2-10
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 5

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 4

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 5

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 4

$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 3

$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 6

$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 7

$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 3
